$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status (column E) and corrector name (column F) for row 2 (index.php entry)
$ws.Range("E2").Value = "P"
$ws.Range("F2").Value = "Krists"

# Align the whole "Labojumu veicējs" column body (F2:F11) like the rest of the table
# (copy the already-centered format from E2 so we reuse the existing style instead
# of minting a new one)
$ws.Range("E2").Copy()
$ws.Range("F2:F11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the selection/view to E2 (also clears the old scrolled-down top-left cell)
$ws.Range("E2").Select()
